$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Start from a clean slate: drop the old hyperlinks and any per-cell
#    formatting so stale styles (e.g. the old text-number-format on row 3)
#    don't leak into the rebuilt grid below.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Cells.ClearFormats()

# ---------------------------------------------------------------------------
# 2) Rewrite the table. New layout adds a "Wait Time" column (G) and two new
#    rows (wait, compare), and renumbers/reorders the Step column.
#    Row 1 = header, rows 2-7 = steps.
# ---------------------------------------------------------------------------

# Header row
$ws.Range("A1").Value = "Step"
$ws.Range("B1").Value = "Accion"
$ws.Range("C1").Value = "Valor Accion"
$ws.Range("D1").Value = "Locator"
$ws.Range("E1").Value = "Valor Locator"
$ws.Range("F1").Value = "Screenshot"
$ws.Range("G1").Value = "Wait Time"

# Row 2: navigate
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "navigate"
$ws.Range("C2").Value = "https://www.google.com/"
$ws.Range("D2").Value = $null
$ws.Range("E2").Value = $null
$ws.Range("F2").Value = $true
$ws.Range("G2").Value = $null

# Row 3: wait
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "wait "
$ws.Range("C3").Value = $null
$ws.Range("D3").Value = $null
$ws.Range("E3").Value = $null
$ws.Range("F3").Value = $true
$ws.Range("G3").Value = 10

# Row 4: type
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "type"
$ws.Range("C4").Value = "Selenium"
$ws.Range("D4").Value = "name"
$ws.Range("E4").Value = "q"
$ws.Range("F4").Value = $true
$ws.Range("G4").Value = $null

# Row 5: click
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "click"
$ws.Range("C5").Value = $null
$ws.Range("D5").Value = "name"
$ws.Range("E5").Value = "btnK"
$ws.Range("F5").Value = $true
$ws.Range("G5").Value = $null

# Row 6: compare
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "compare "
$ws.Range("C6").Value = "hola"
$ws.Range("D6").Value = "xpath"
$ws.Range("E6").Value = "//div[@class='kno-ecr-pt PZPZlf gsmt i8lZMc']//span[contains(text(),'Selenium')]"
$ws.Range("F6").Value = $true
$ws.Range("G6").Value = $null

# Row 7: quit
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "quit"
$ws.Range("C7").Value = $null
$ws.Range("D7").Value = $null
$ws.Range("E7").Value = $null
$ws.Range("F7").Value = $true
$ws.Range("G7").Value = $null

# ---------------------------------------------------------------------------
# 3) Hyperlink for C2 (navigate -> google).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C2"), "https://www.google.com/") | Out-Null

# ---------------------------------------------------------------------------
# 4) Formatting: border around the whole table, accent fill for header row +
#    step column, text number format on header + the "type" row values.
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:G7")
$tableRange.Borders.LineStyle = 1
$tableRange.Borders.Weight = 2

$ws.Range("A1:G1").NumberFormat = "@"
$ws.Range("A1:G1").Interior.ThemeColor = 5

$ws.Range("A2:A7").Interior.ThemeColor = 5

$ws.Range("B4:E4").NumberFormat = "@"

# E6 keeps the default (no border) style, matching the source edit.
$ws.Range("E6").Borders.LineStyle = -4142

# ---------------------------------------------------------------------------
# 5) Selection / active cell.
# ---------------------------------------------------------------------------
$ws.Range("E8").Select()
